$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.869.95'
$ws.Range("E2").Value = '  +0.34%  '

# Row 3
$ws.Range("D3").Value = '2.288.97'
$ws.Range("E3").Value = '  -0.03%  '

# Row 4
$ws.Range("E4").Value = '  +0.35%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '113.55'
$ws.Range("E5").Value = '  +17.81%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '269.23'
$ws.Range("E6").Value = '  +0.02%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("E7").Value = '  +0.24%  '

# Row 8
$ws.Range("E8").Value = '  +0.36%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.617'
$ws.Range("E9").Value = '  +1.33%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.61'
$ws.Range("E10").Value = '  +4.23%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0945'
$ws.Range("E11").Value = '  +0.79%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.99'
$ws.Range("E12").Value = '  +13.07%  '

# Row 13
$ws.Range("E13").Value = '  -0.07%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.78'
$ws.Range("E14").Value = '  +0.81%  '

# Row 15
$ws.Range("D15").Value = '2.638.95'
$ws.Range("E15").Value = '  +0.20%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.848'
$ws.Range("E16").Value = '  -0.91%  '

# Row 17
$ws.Range("D17").Value = '2.293.77'
$ws.Range("E17").Value = '  +0.21%  '

# Row 18
$ws.Range("D18").Value = '43.832.58'
$ws.Range("E18").Value = '  +0.38%  '

# Row 19
$ws.Range("E19").Value = '  -2.02%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.78'
$ws.Range("E20").Value = '  +9.66%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.36'
$ws.Range("E21").Value = '  +0.30%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.43'
$ws.Range("E22").Value = '  -1.38%  '

# Row 23
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.30'
$ws.Range("E23").Value = '  -0.23%  '

# Row 24
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.96'
$ws.Range("E24").Value = '  +10.13%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.64'
$ws.Range("E25").Value = '  +6.28%  '

# Row 26
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.02%  '

# Row 27
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.60'
$ws.Range("E27").Value = '  +2.49%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '41.45'
$ws.Range("E28").Value = '  +7.18%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.39'
$ws.Range("E29").Value = '  -1.70%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.26'
$ws.Range("E30").Value = '  +1.50%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.11'
$ws.Range("E31").Value = '  +0.21%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0930'
$ws.Range("E32").Value = '  +3.80%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.52'
$ws.Range("E33").Value = '  -3.08%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.67'
$ws.Range("E34").Value = '  +4.28%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.127'
$ws.Range("E35").Value = '  +0.12%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.63'
$ws.Range("E36").Value = '  -0.32%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0361'
$ws.Range("E37").Value = '  +2.69%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.107'
$ws.Range("E38").Value = '  +0.60%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.77'
$ws.Range("E39").Value = '  +5.46%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '73.96'
$ws.Range("E40").Value = '  +14.44%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.243'
$ws.Range("E41").Value = '  +3.02%  '

# Row 42
$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.37'
$ws.Range("E42").Value = '  +2.66%  '

# Row 43
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.50'
$ws.Range("E43").Value = '  +9.26%  '

# Row 44
$ws.Range("B44").Value = 'THORChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.26'
$ws.Range("E44").Value = '  +20.55%  '

# Row 45
$ws.Range("E45").Value = '  +0.28%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.39'
$ws.Range("E46").Value = '  +3.80%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.80'
$ws.Range("E47").Value = '  +0.88%  '

# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0997'
$ws.Range("E48").Value = '  -2.39%  '

# Row 49
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.37'
$ws.Range("E49").Value = '  +3.99%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.466'
$ws.Range("E50").Value = '  +7.81%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.23'
$ws.Range("E51").Value = '  +2.21%  '
